$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.812.48"
$ws.Range("E2").Value = "  -0.60%  "

# Row 3
$ws.Range("D3").Value = "2.924.57"
$ws.Range("E3").Value = "  +1.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.01"
$ws.Range("E5").Value = "  +0.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.65"
$ws.Range("E6").Value = "  -1.18%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  +1.99%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +1.27%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.29"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0884"
$ws.Range("E11").Value = "  +3.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.137"
$ws.Range("E12").Value = "  +0.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.71"
$ws.Range("E13").Value = "  -1.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.89"
$ws.Range("E14").Value = "  +1.29%  "

# Row 15
$ws.Range("D15").Value = "3.384.46"
$ws.Range("E15").Value = "  +1.07%  "

# Row 16
$ws.Range("D16").Value = "2.906.37"
$ws.Range("E16").Value = "  +1.46%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.981"
$ws.Range("E17").Value = "  -1.74%  "

# Row 18
$ws.Range("D18").Value = "51.809.31"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("E19").Value = "  -1.41%  "

# Row 20
$ws.Range("E20").Value = "  -2.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.03"
$ws.Range("E21").Value = "  -2.07%  "

# Row 22
$ws.Range("E22").Value = "  +0.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.83"
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.57"
$ws.Range("E24").Value = "  +0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  +1.04%  "

# Row 26
$ws.Range("E26").Value = "  +10.60%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.12"
$ws.Range("E27").Value = "  +2.37%  "

# Row 29
$ws.Range("E29").Value = "  +12.93%  "

# Row 30
$ws.Range("E30").Value = "  +11.80%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.56"
$ws.Range("E31").Value = "  +0.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.72"
$ws.Range("E32").Value = "  -0.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.07"
$ws.Range("E33").Value = "  -1.50%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.20"
$ws.Range("E34").Value = "  -1.92%  "

# Row 35
$ws.Range("E35").Value = "  -4.41%  "

# Row 36
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("E37").Value = "  -15.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  -2.80%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.40"
$ws.Range("E39").Value = "  -1.56%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("E40").Value = "  +6.05%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.00"
$ws.Range("E41").Value = "  -2.18%  "

# Row 42
$ws.Range("E42").Value = "  +2.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.03"
$ws.Range("E43").Value = "  +1.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.07"
$ws.Range("E44").Value = "  -1.71%  "

# Row 45
$ws.Range("E45").Value = "  -1.33%  "

# Row 46
$ws.Range("E46").Value = "  +0.27%  "

# Row 47
$ws.Range("E47").Value = "  -4.05%  "

# Row 48
$ws.Range("D48").Value = "2.139.12"
$ws.Range("E48").Value = "  -3.03%  "

# Row 49
$ws.Range("E49").Value = "  -5.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0338"
$ws.Range("E50").Value = "  +5.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.906"
$ws.Range("E51").Value = "  -4.40%  "
